# Applies the 2023-11-11 update to the Poland Division 2 2023-2024 sheet:
#  - Swap the match data (columns F:V) between rows 92 and 93
#    (two games that had been entered against the wrong row).
#  - Swap the match data (columns F:V) between rows 136 and 137
#    (same kind of correction for another pair of matches).
#  - Append a new match result as row 145 (Ol. Grudziadz vs Hutnik Krakow).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns F (6) through V (22) hold the match-specific data; A:E (1:5) are
# the row index / pais / torneio / temporada / data_partida, which stay put.
for ($col = 6; $col -le 22; $col++) {
    $cellA = $ws.Cells.Item(92, $col)
    $cellB = $ws.Cells.Item(93, $col)
    $valA = $cellA.Value()
    $valB = $cellB.Value()
    $cellA.Value = $valB
    $cellB.Value = $valA
}

for ($col = 6; $col -le 22; $col++) {
    $cellA = $ws.Cells.Item(136, $col)
    $cellB = $ws.Cells.Item(137, $col)
    $valA = $cellA.Value()
    $valB = $cellB.Value()
    $cellA.Value = $valB
    $cellB.Value = $valA
}

# Append the new match as row 145. First clone the formatting used by the
# other data rows (bold/bordered index cell in column A, datetime format in
# column E) from the previous row, then fill in the values.
$ws.Cells.Item(144, 1).Copy()
$ws.Cells.Item(145, 1).PasteSpecial(-4122)
$ws.Cells.Item(144, 5).Copy()
$ws.Cells.Item(145, 5).PasteSpecial(-4122)

$newRow = 145
$ws.Cells.Item($newRow, 1).Value = 144
$ws.Cells.Item($newRow, 2).Value = "poland"
$ws.Cells.Item($newRow, 3).Value = "division-2"
$ws.Cells.Item($newRow, 4).Value = "2023-2024"
$ws.Cells.Item($newRow, 5).Value = 45241.54166666666
$ws.Cells.Item($newRow, 6).Value = "Ol. Grudziadz"
$ws.Cells.Item($newRow, 7).Value = 0
$ws.Cells.Item($newRow, 8).Value = "Hutnik Krakow"
$ws.Cells.Item($newRow, 9).Value = 2
$ws.Cells.Item($newRow, 10).Value = 2.09
$ws.Cells.Item($newRow, 11).Value = "10/11/2023 01:13"
$ws.Cells.Item($newRow, 12).Value = 1.74
$ws.Cells.Item($newRow, 13).Value = "11/11/2023 12:51"
$ws.Cells.Item($newRow, 14).Value = 3.28
$ws.Cells.Item($newRow, 15).Value = "10/11/2023 01:13"
$ws.Cells.Item($newRow, 16).Value = 3.86
$ws.Cells.Item($newRow, 17).Value = "11/11/2023 12:51"
$ws.Cells.Item($newRow, 18).Value = 3.02
$ws.Cells.Item($newRow, 19).Value = "10/11/2023 01:13"
$ws.Cells.Item($newRow, 20).Value = 4.2
$ws.Cells.Item($newRow, 21).Value = "11/11/2023 12:51"
$ws.Cells.Item($newRow, 22).Value = "https://www.betexplorer.com/football/poland/division-2/ol-grudziadz-hutnik-krakow/OQJpzDR2/"
